# Update "想去人数" (F column) counts for rows 2-21 on the "展览" and "全部类型" sheets.
# Mapping of row -> new value (only rows whose value actually changed).
$updates = @{
    2  = 302
    3  = 309
    4  = 58
    5  = 367
    6  = 11078
    7  = 540
    8  = 97
    10 = 83
    12 = 148
    13 = 18
    14 = 47
    15 = 43
    17 = 32
    18 = 312
    19 = 1176
    20 = 57
    21 = 884
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
